# Osprey EMEA gold workbook - "product qty for osprey"
#
# On the "Review" sheet, insert a new "MoreQuantity" column right before the
# "FirstName" column (column G), add its header, and populate the product
# quantity values for the configurable-product row (row 4): "000" in the new
# MoreQuantity column and "10+" in the (now shifted) Quantity column.
#
# Both values are entered the way a person typing into Excel would force
# text (leading apostrophe) so they keep their literal digits/"+" instead of
# being coerced to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review")

# Shift FirstName..Review (G:Z) one column to the right, creating a blank
# column G for the new field.
$ws.Columns("G").Insert()

# New column header.
$ws.Range("G1").Value = "MoreQuantity"

# Quantity value for the AETHER(TM) 55 configurable product row moved from
# column S to column T after the insert.
$ws.Range("T4").Value = "'10+"

# New MoreQuantity value for the same row.
$ws.Range("G4").Value = "'000"

# Make "Review" the active sheet with the newly edited cell selected, as it
# would be right after making this edit.
$ws.Activate()
$ws.Range("G4").Select()
